$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = -20.575
$ws.Range("C4").Value = -12.566
$ws.Range("A6").Value = -22.157
$ws.Range("A7").Value = -19.81
$ws.Range("C9").Value = -11.361
$ws.Range("C12").Value = -10.919
$ws.Range("A16").Value = -22.046
$ws.Range("C17").Value = -13.261
$ws.Range("C18").Value = -11.591
$ws.Range("C19").Value = -12.126
$ws.Range("A20").Value = -19.828
$ws.Range("C20").Value = -11.729
$ws.Range("C26").Value = -12.222
$ws.Range("A28").Value = -21.934
$ws.Range("A29").Value = -21.344
$ws.Range("C31").Value = -13.298
$ws.Range("A32").Value = -21.772
$ws.Range("C39").Value = -12.133
$ws.Range("A40").Value = -20.007
$ws.Range("C40").Value = -12.226
$ws.Range("C41").Value = -12.117
$ws.Range("C42").Value = -12.29
$ws.Range("C43").Value = -12.566
$ws.Range("A46").Value = -21.89
$ws.Range("C47").Value = -12.082
$ws.Range("C48").Value = -11.74
$ws.Range("A51").Value = -21.715
$ws.Range("A52").Value = -21.972
$ws.Range("A57").Value = -22.202
$ws.Range("A59").Value = -22.439
$ws.Range("A62").Value = -21.95
$ws.Range("C63").Value = -11.178
$ws.Range("C64").Value = -11.011
$ws.Range("A66").Value = -21.567
$ws.Range("A73").Value = -20.345
$ws.Range("A74").Value = -21.197
$ws.Range("C76").Value = -12.969
$ws.Range("C81").Value = -13.176
$ws.Range("C89").Value = -13.626
$ws.Range("A92").Value = -21.759
$ws.Range("C94").Value = -11.539
$ws.Range("A100").Value = -22.143
